$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H80").Value = 809.63635
$ws.Range("I80").Value = 567.1667
$ws.Range("J80").Value = 1100.6
$ws.Range("K80").Value = 1701.5001
$ws.Range("L80").Value = 3301.8
$ws.Range("M80").Value = -703.5001
$ws.Range("N80").Value = -5297.799999999999

$ws.Range("H83").Value = 809.63635
$ws.Range("I83").Value = 567.1667
$ws.Range("J83").Value = 1100.6
$ws.Range("K83").Value = 5104.5003
$ws.Range("L83").Value = 9905.4
$ws.Range("M83").Value = -112.5002999999997
$ws.Range("N83").Value = -19889.4

$ws.Range("H100").Value = 71429870
$ws.Range("I100").Value = 1517.5
$ws.Range("K100").Value = 1517.5
$ws.Range("M100").Value = -976.5

$ws.Range("H125").Value = 1818.2858
$ws.Range("I125").Value = 516.6667
$ws.Range("J125").Value = 2173.2727
$ws.Range("K125").Value = 4650.0003
$ws.Range("L125").Value = 19559.4543
$ws.Range("M125").Value = -2190.0003
$ws.Range("N125").Value = -24479.4543

$ws.Range("H126").Value = 41000
$ws.Range("J126").Value = 41000
$ws.Range("L126").Value = 41000
$ws.Range("N126").Value = -50880

$ws.Range("H128").Value = 34494
$ws.Range("J128").Value = 34494
$ws.Range("L128").Value = 34494
$ws.Range("N128").Value = -44454

$ws.Range("H141").Value = 3283.3333
$ws.Range("I141").Value = 3425
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 10275
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -5095
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets("ARM")
$ws.Range("H61").Value = 1542.375
$ws.Range("I61").Value = 1336.7693
$ws.Range("J61").Value = 2433.3333
$ws.Range("K61").Value = 1336.7693
$ws.Range("L61").Value = 2433.3333
$ws.Range("M61").Value = -1124.7693
$ws.Range("N61").Value = -2857.3333

$ws.Range("H103").Value = 35371.668
$ws.Range("J103").Value = 35371.668
$ws.Range("L103").Value = 35371.668
$ws.Range("N103").Value = -37715.668

$ws.Range("H109").Value = 16919.25
$ws.Range("J109").Value = 16919.25
$ws.Range("L109").Value = 16919.25
$ws.Range("N109").Value = -19693.25

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H132").Value = 2444.76
$ws.Range("I132").Value = 2271.5625
$ws.Range("J132").Value = 2752.6667
$ws.Range("K132").Value = 6814.6875
$ws.Range("L132").Value = 8258.000100000001
$ws.Range("M132").Value = -4284.6875
$ws.Range("N132").Value = -13318.0001

$ws.Range("H136").Value = 1542.375
$ws.Range("I136").Value = 1336.7693
$ws.Range("J136").Value = 2433.3333
$ws.Range("K136").Value = 4010.3079
$ws.Range("L136").Value = 7299.999899999999
$ws.Range("M136").Value = -1460.3079
$ws.Range("N136").Value = -12399.9999

$ws = $wb.Worksheets("BSM")
$ws.Range("H134").Value = 1411.037
$ws.Range("I134").Value = 1373.3684
$ws.Range("J134").Value = 1500.5
$ws.Range("K134").Value = 4120.1052
$ws.Range("L134").Value = 4501.5
$ws.Range("M134").Value = -1585.1052
$ws.Range("N134").Value = -9571.5

$ws = $wb.Worksheets("CUL")
$ws.Range("H33").Value = 483.23077
$ws.Range("I33").Value = 432.85715
$ws.Range("J33").Value = 542
$ws.Range("K33").Value = 2597.1429
$ws.Range("L33").Value = 3252
$ws.Range("M33").Value = -2314.1429
$ws.Range("N33").Value = -3818

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H44").Value = 645.16364
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 645.16364
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 1935.49092
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -2731.49092

$ws.Range("H68").Value = 3480.6
$ws.Range("I68").Value = 5300.6665
$ws.Range("K68").Value = 15901.9995
$ws.Range("M68").Value = -15090.9995

$ws.Range("H71").Value = 3480.6
$ws.Range("I71").Value = 5300.6665
$ws.Range("K71").Value = 47705.9985
$ws.Range("M71").Value = -43649.9985

$ws = $wb.Worksheets("GSM")
$ws.Range("H122").Value = 1963.125
$ws.Range("I122").Value = 2126.75
$ws.Range("J122").Value = 1799.5
$ws.Range("K122").Value = 6380.25
$ws.Range("L122").Value = 5398.5
$ws.Range("M122").Value = -3930.25
$ws.Range("N122").Value = -10298.5

$ws.Range("H124").Value = 55997.5
$ws.Range("J124").Value = 55997.5
$ws.Range("L124").Value = 55997.5
$ws.Range("N124").Value = -65817.5

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 5557577.5
$ws.Range("I7").Value = 9092712
$ws.Range("J7").Value = 2366.2856
$ws.Range("K7").Value = 9092712
$ws.Range("L7").Value = 2366.2856
$ws.Range("M7").Value = -9092600
$ws.Range("N7").Value = -2590.2856

$ws.Range("H40").Value = 1910.3
$ws.Range("I40").Value = 1904.8572
$ws.Range("J40").Value = 1923
$ws.Range("K40").Value = 1904.8572
$ws.Range("L40").Value = 1923
$ws.Range("M40").Value = -1768.8572
$ws.Range("N40").Value = -2195

$ws.Range("H100").Value = 1742.8572
$ws.Range("I100").Value = 1742.8572
$ws.Range("K100").Value = 1742.8572
$ws.Range("M100").Value = -1201.8572

$ws.Range("H108").Value = 20350
$ws.Range("J108").Value = 20350
$ws.Range("L108").Value = 20350
$ws.Range("N108").Value = -28030

$ws.Range("H122").Value = 4230.7334
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 4663.4165
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 13990.2495
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -18890.2495

$ws.Range("H126").Value = 5557577.5
$ws.Range("I126").Value = 9092712
$ws.Range("J126").Value = 2366.2856
$ws.Range("K126").Value = 27278136
$ws.Range("L126").Value = 7098.8568
$ws.Range("M126").Value = -27275666
$ws.Range("N126").Value = -12038.8568

$ws.Range("H127").Value = 56428.332
$ws.Range("J127").Value = 56428.332
$ws.Range("L127").Value = 56428.332
$ws.Range("N127").Value = -66348.33199999999

$ws = $wb.Worksheets("WVR")
$ws.Range("H100").Value = 770674.5600000001
$ws.Range("I100").Value = 1701.6
$ws.Range("J100").Value = 3333917.8
$ws.Range("K100").Value = 3403.2
$ws.Range("L100").Value = 6667835.6
$ws.Range("M100").Value = -2862.2
$ws.Range("N100").Value = -6668917.6

$ws.Range("H135").Value = 67216.42999999999
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 67216.42999999999
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 67216.42999999999
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -77356.42999999999

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()
